$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 26: update title only
$ws.Range("D26").Value = "ai plus(est soft)"

# Row 36: update title and link
$ws.Range("D36").Value = "Self-Supervised Learning Methods for Chemical Property Prediction"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/361"

# Row 46: update title and link
$ws.Range("D46").Value = "[Bioinformatics] 2021년 12월,  유전체 정보분석 전문가 기본과정 II"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/458"

# Row 50: update title only
$ws.Range("D50").Value = "가장 쉬운 Mahalanobis distance"
